$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.233789666666667
$ws.Range("H2").Value = 3.701369
$ws.Range("I2").Value = 0.0001664233864291757
$ws.Range("J2").Value = 0.0001664233864291757
$ws.Range("M2").Value = 1.949849666666667
$ws.Range("N2").Value = 5.849549000000001
$ws.Range("O2").Value = 0.06676506732104066
$ws.Range("P2").Value = 0.06676506732104066
$ws.Range("Q2").Value = 2.405704370286778
$ws.Range("R2").Value = 21.651339332581
$ws.Range("S2").Value = 0.00001111126859873948
$ws.Range("T2").Value = 0.00001111126859873948

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.233789666666667
$ws.Range("H3").Value = 3.701369
$ws.Range("I3").Value = 0.0001664233864291757
$ws.Range("J3").Value = 0.0001664233864291757
$ws.Range("O3").Value = 0.7967262871802238
$ws.Range("P3").Value = 0.7967262871802239
$ws.Range("Q3").Value = 28.70794545559889
$ws.Range("R3").Value = 258.37150910039
$ws.Range("S3").Value = 0.0001325938867696768
$ws.Range("T3").Value = 0.0001325938867696768

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.233789666666667
$ws.Range("H4").Value = 3.701369
$ws.Range("I4").Value = 0.0001664233864291757
$ws.Range("J4").Value = 0.0001664233864291757
$ws.Range("O4").Value = 0.1365086454987356
$ws.Range("P4").Value = 0.1365086454987356
$ws.Range("Q4").Value = 4.918731579781444
$ws.Range("R4").Value = 44.268584218033
$ws.Range("S4").Value = 0.00002271823106075944
$ws.Range("T4").Value = 0.00002271823106075944

$ws.Range("I5").Value = 0.9827534361704352
$ws.Range("J5").Value = 0.9827534361704352
$ws.Range("M5").Value = 1.949849666666667
$ws.Range("N5").Value = 5.849549000000001
$ws.Range("O5").Value = 0.06676506732104066
$ws.Range("P5").Value = 0.06676506732104066
$ws.Range("Q5").Value = 14206.02168383165
$ws.Range("R5").Value = 127854.1951544848
$ws.Range("S5").Value = 0.06561359932590315
$ws.Range("T5").Value = 0.06561359932590315

$ws.Range("I6").Value = 0.9827534361704352
$ws.Range("J6").Value = 0.9827534361704352
$ws.Range("O6").Value = 0.7967262871802238
$ws.Range("P6").Value = 0.7967262871802239
$ws.Range("S6").Value = 0.7829854964136779
$ws.Range("T6").Value = 0.782985496413678

$ws.Range("I7").Value = 0.9827534361704352
$ws.Range("J7").Value = 0.9827534361704352
$ws.Range("O7").Value = 0.1365086454987356
$ws.Range("P7").Value = 0.1365086454987356
$ws.Range("S7").Value = 0.1341543404308543
$ws.Range("T7").Value = 0.1341543404308543

$ws.Range("I8").Value = 0.01708014044313564
$ws.Range("J8").Value = 0.01708014044313564
$ws.Range("M8").Value = 1.949849666666667
$ws.Range("N8").Value = 5.849549000000001
$ws.Range("O8").Value = 0.06676506732104066
$ws.Range("P8").Value = 0.06676506732104066
$ws.Range("Q8").Value = 246.8990049463379
$ws.Range("R8").Value = 2222.091044517041
$ws.Range("S8").Value = 0.00114035672653878
$ws.Range("T8").Value = 0.00114035672653878

$ws.Range("I9").Value = 0.01708014044313564
$ws.Range("J9").Value = 0.01708014044313564
$ws.Range("O9").Value = 0.7967262871802238
$ws.Range("P9").Value = 0.7967262871802239
$ws.Range("S9").Value = 0.01360819687977624
$ws.Range("T9").Value = 0.01360819687977624

$ws.Range("I10").Value = 0.01708014044313564
$ws.Range("J10").Value = 0.01708014044313564
$ws.Range("O10").Value = 0.1365086454987356
$ws.Range("P10").Value = 0.1365086454987356
$ws.Range("S10").Value = 0.002331586836820621
$ws.Range("T10").Value = 0.002331586836820621

